$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, bordered, centered) from an existing
# header cell (AC1) onto the three new header cells, then set their text.
$headerFormat = $ws.Range("AC1")
$headerFormat.Copy()

$ws.Range("AD1").PasteSpecial(-4122)
$ws.Range("AE1").PasteSpecial(-4122)
$ws.Range("AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record (Wins/Losses/Ties) is the same for every player row in
# this sheet, so fan the three values out across all data rows.
for ($row = 2; $row -le 53; $row++) {
    $ws.Cells.Item($row, 30).Value = 85
    $ws.Cells.Item($row, 31).Value = 77
    $ws.Cells.Item($row, 32).Value = 0
}

Write-Output "done"
